$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 13913.1
$ws.Range("I32").Value = 7000.3335
$ws.Range("K32").Value = 7000.3335
$ws.Range("M32").Value = -6674.3335
$ws.Range("H33").Value = 519
$ws.Range("I33").Value = 519
$ws.Range("K33").Value = 519
$ws.Range("M33").Value = -290
$ws.Range("H38").Value = 1105.5
$ws.Range("I38").Value = 263.57144
$ws.Range("J38").Value = 6999
$ws.Range("K38").Value = 790.71432
$ws.Range("L38").Value = 20997
$ws.Range("M38").Value = -418.71432
$ws.Range("N38").Value = -21741
$ws.Range("H76").Value = 37044556
$ws.Range("I76").Value = 55563840
$ws.Range("K76").Value = 55563840
$ws.Range("M76").Value = -55563525
$ws.Range("H79").Value = 37044556
$ws.Range("I79").Value = 55563840
$ws.Range("K79").Value = 55563840
$ws.Range("M79").Value = -55562748
$ws.Range("H80").Value = 1574
$ws.Range("I80").Value = 1574
$ws.Range("K80").Value = 4722
$ws.Range("M80").Value = -3724
$ws.Range("H83").Value = 1574
$ws.Range("I83").Value = 1574
$ws.Range("K83").Value = 14166
$ws.Range("M83").Value = -9174
$ws.Range("H88").Value = 3860.2222
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 5668.4
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 5668.4
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -6480.4
$ws.Range("H91").Value = 3860.2222
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 5668.4
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 5668.4
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -8476.4
$ws.Range("I96").Value = 1190.55
$ws.Range("K96").Value = 3571.65
$ws.Range("M96").Value = -2198.65
$ws.Range("H106").Value = 6051
$ws.Range("I106").Value = 5451.6665
$ws.Range("J106").Value = 6650.3335
$ws.Range("K106").Value = 5451.6665
$ws.Range("L106").Value = 6650.3335
$ws.Range("M106").Value = -4820.6665
$ws.Range("N106").Value = -7912.3335
$ws.Range("H107").Value = 1387.1578
$ws.Range("I107").Value = 1387.1578
$ws.Range("K107").Value = 1387.1578
$ws.Range("M107").Value = 532.8422
$ws.Range("H115").Value = 1994.1666
$ws.Range("J115").Value = 9900
$ws.Range("L115").Value = 29700
$ws.Range("N115").Value = -32834
$ws.Range("H117").Value = 70555
$ws.Range("I117").Value = 70555
$ws.Range("K117").Value = 70555
$ws.Range("M117").Value = -65966
$ws.Range("H118").Value = 604
$ws.Range("H127").Value = 21171.166
$ws.Range("I127").Value = 17998.5
$ws.Range("K127").Value = 53995.5
$ws.Range("M127").Value = -49035.5
$ws.Range("H132").Value = 317901.1
$ws.Range("I132").Value = 350339.28
$ws.Range("J132").Value = 1628.75
$ws.Range("K132").Value = 1051017.84
$ws.Range("L132").Value = 4886.25
$ws.Range("M132").Value = -1048487.84
$ws.Range("N132").Value = -9946.25
$ws.Range("H137").Value = 7956.1
$ws.Range("I137").Value = 9084.909
$ws.Range("J137").Value = 6576.4443
$ws.Range("K137").Value = 27254.727
$ws.Range("L137").Value = 19729.3329
$ws.Range("M137").Value = -24704.727
$ws.Range("N137").Value = -24829.3329
$ws.Range("H139").Value = 117259.836
$ws.Range("J139").Value = 114519.664
$ws.Range("L139").Value = 114519.664
$ws.Range("N139").Value = -124799.664
$ws.Range("H141").Value = 3956.1428
$ws.Range("J141").Value = 5564.6665
$ws.Range("L141").Value = 16693.9995
$ws.Range("N141").Value = -27053.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3572.0557
$ws.Range("I2").Value = 1674.8125
$ws.Range("K2").Value = 1674.8125
$ws.Range("M2").Value = -1561.8125
$ws.Range("H32").Value = 9260231
$ws.Range("I32").Value = 7693295
$ws.Range("J32").Value = 111111110
$ws.Range("K32").Value = 7693295
$ws.Range("L32").Value = 111111110
$ws.Range("M32").Value = -7693008
$ws.Range("N32").Value = -111111684
$ws.Range("H34").Value = 40000
$ws.Range("I34").Value = 50000
$ws.Range("J34").Value = 30000
$ws.Range("K34").Value = 50000
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = -49729
$ws.Range("N34").Value = -30542
$ws.Range("H45").Value = 1334.2
$ws.Range("I45").Value = 1305.25
$ws.Range("K45").Value = 1305.25
$ws.Range("M45").Value = -928.25
$ws.Range("H61").Value = 5206.8037
$ws.Range("I61").Value = 6756.9165
$ws.Range("K61").Value = 6756.9165
$ws.Range("M61").Value = -6544.9165
$ws.Range("H74").Value = 7095
$ws.Range("I74").Value = 7588.8887
$ws.Range("K74").Value = 7588.8887
$ws.Range("M74").Value = -6714.8887
$ws.Range("H77").Value = 7095
$ws.Range("I77").Value = 7588.8887
$ws.Range("K77").Value = 37944.4435
$ws.Range("M77").Value = -33576.4435
$ws.Range("H102").Value = 2398.1428
$ws.Range("I102").Value = 2214.7778
$ws.Range("K102").Value = 2214.7778
$ws.Range("M102").Value = -592.7777999999998
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 5662.4585
$ws.Range("I110").Value = 2756.5
$ws.Range("J110").Value = 7115.4375
$ws.Range("K110").Value = 2756.5
$ws.Range("L110").Value = 7115.4375
$ws.Range("M110").Value = -711.5
$ws.Range("N110").Value = -11205.4375
$ws.Range("H116").Value = 3572.0557
$ws.Range("I116").Value = 1674.8125
$ws.Range("K116").Value = 1674.8125
$ws.Range("M116").Value = 619.1875
$ws.Range("H132").Value = 630219.9
$ws.Range("I132").Value = 659773.5600000001
$ws.Range("K132").Value = 1979320.68
$ws.Range("M132").Value = -1976790.68
$ws.Range("H136").Value = 5206.8037
$ws.Range("I136").Value = 6756.9165
$ws.Range("K136").Value = 20270.7495
$ws.Range("M136").Value = -17720.7495

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3572.0557
$ws.Range("I3").Value = 1674.8125
$ws.Range("K3").Value = 1674.8125
$ws.Range("M3").Value = -1560.8125
$ws.Range("H64").Value = 4257.5
$ws.Range("I64").Value = 1320.8334
$ws.Range("J64").Value = 13067.5
$ws.Range("K64").Value = 1320.8334
$ws.Range("L64").Value = 13067.5
$ws.Range("M64").Value = -1095.8334
$ws.Range("N64").Value = -13517.5
$ws.Range("H67").Value = 4257.5
$ws.Range("I67").Value = 1320.8334
$ws.Range("J67").Value = 13067.5
$ws.Range("K67").Value = 1320.8334
$ws.Range("L67").Value = 13067.5
$ws.Range("M67").Value = -540.8334
$ws.Range("N67").Value = -14627.5
$ws.Range("H86").Value = 3370.9583
$ws.Range("I86").Value = 1207.7222
$ws.Range("J86").Value = 9860.666999999999
$ws.Range("K86").Value = 1207.7222
$ws.Range("L86").Value = 9860.666999999999
$ws.Range("M86").Value = -84.72219999999993
$ws.Range("N86").Value = -12106.667
$ws.Range("H89").Value = 3370.9583
$ws.Range("I89").Value = 1207.7222
$ws.Range("J89").Value = 9860.666999999999
$ws.Range("K89").Value = 6038.611
$ws.Range("L89").Value = 49303.335
$ws.Range("M89").Value = -422.6109999999999
$ws.Range("N89").Value = -60535.335
$ws.Range("H99").Value = 5091.921
$ws.Range("J99").Value = 8930.117
$ws.Range("L99").Value = 8930.117
$ws.Range("N99").Value = -11926.117
$ws.Range("H105").Value = 2563.0303
$ws.Range("I105").Value = 2793.4211
$ws.Range("K105").Value = 2793.4211
$ws.Range("M105").Value = -1046.4211
$ws.Range("H107").Value = 2225785.5
$ws.Range("I107").Value = 2943716.5
$ws.Range("J107").Value = 6726.273
$ws.Range("K107").Value = 2943716.5
$ws.Range("L107").Value = 6726.273
$ws.Range("M107").Value = -2941796.5
$ws.Range("N107").Value = -10566.273
$ws.Range("H134").Value = 930715.3
$ws.Range("I134").Value = 965406.25
$ws.Range("K134").Value = 2896218.75
$ws.Range("M134").Value = -2893683.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 16670507
$ws.Range("I16").Value = 45456640
$ws.Range("K16").Value = 45456640
$ws.Range("M16").Value = -45456353
$ws.Range("H31").Value = 4477.3887
$ws.Range("I31").Value = 1189.2858
$ws.Range("J31").Value = 5271.069
$ws.Range("K31").Value = 1189.2858
$ws.Range("L31").Value = 5271.069
$ws.Range("M31").Value = -894.2858000000001
$ws.Range("N31").Value = -5861.069
$ws.Range("H34").Value = 4477.3887
$ws.Range("I34").Value = 1189.2858
$ws.Range("J34").Value = 5271.069
$ws.Range("K34").Value = 1189.2858
$ws.Range("L34").Value = 5271.069
$ws.Range("M34").Value = -987.2858000000001
$ws.Range("N34").Value = -5675.069
$ws.Range("H58").Value = 18877158
$ws.Range("I58").Value = 37042540
$ws.Range("J58").Value = 13108.27
$ws.Range("K58").Value = 37042540
$ws.Range("L58").Value = 13108.27
$ws.Range("M58").Value = -37042337
$ws.Range("N58").Value = -13514.27
$ws.Range("H62").Value = 6339.6
$ws.Range("I62").Value = 6258.636
$ws.Range("K62").Value = 6258.636
$ws.Range("M62").Value = -5634.636
$ws.Range("H65").Value = 6339.6
$ws.Range("I65").Value = 6258.636
$ws.Range("K65").Value = 31293.18
$ws.Range("M65").Value = -28173.18
$ws.Range("H74").Value = 25985
$ws.Range("J74").Value = 25985
$ws.Range("L74").Value = 25985
$ws.Range("N74").Value = -27733
$ws.Range("H77").Value = 25985
$ws.Range("J77").Value = 25985
$ws.Range("L77").Value = 77955
$ws.Range("N77").Value = -86691
$ws.Range("H94").Value = 43481484
$ws.Range("I94").Value = 76924760
$ws.Range("K94").Value = 76924760
$ws.Range("M94").Value = -76924309
$ws.Range("H107").Value = 1138.28
$ws.Range("J107").Value = 1547.9
$ws.Range("L107").Value = 1547.9
$ws.Range("N107").Value = -5387.9
$ws.Range("H113").Value = 16670507
$ws.Range("I113").Value = 45456640
$ws.Range("K113").Value = 45456640
$ws.Range("M113").Value = -45454470
$ws.Range("H132").Value = 9330.233
$ws.Range("I132").Value = 5593.0386
$ws.Range("J132").Value = 33622
$ws.Range("K132").Value = 16779.1158
$ws.Range("L132").Value = 100866
$ws.Range("M132").Value = -14249.1158
$ws.Range("N132").Value = -105926
$ws.Range("H134").Value = 71434820
$ws.Range("I134").Value = 76927510
$ws.Range("K134").Value = 230782530
$ws.Range("M134").Value = -230779995
$ws.Range("H136").Value = 18877158
$ws.Range("I136").Value = 37042540
$ws.Range("J136").Value = 13108.27
$ws.Range("K136").Value = 111127620
$ws.Range("L136").Value = 39324.81
$ws.Range("M136").Value = -111125070
$ws.Range("N136").Value = -44424.81

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1006.3125
$ws.Range("J5").Value = 1190.1
$ws.Range("L5").Value = 3570.3
$ws.Range("N5").Value = -3794.3
$ws.Range("H35").Value = 240
$ws.Range("I35").Value = 240
$ws.Range("K35").Value = 720
$ws.Range("M35").Value = -432
$ws.Range("H49").Value = 2060
$ws.Range("I49").Value = 2060
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 6180
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -6024
$ws.Range("N49").ClearContents()
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 90000
$ws.Range("N106").Value = -91892
$ws.Range("H120").Value = 37750
$ws.Range("I120").Value = 1000
$ws.Range("K120").Value = 3000
$ws.Range("M120").Value = 1838
$ws.Range("H122").Value = 5971.6313
$ws.Range("I122").Value = 999.6667
$ws.Range("K122").Value = 8997.0003
$ws.Range("M122").Value = -6547.0003
$ws.Range("H134").Value = 7806.115
$ws.Range("I134").Value = 7318.36
$ws.Range("K134").Value = 21955.08
$ws.Range("M134").Value = -16885.08
$ws.Range("H135").Value = 1006.3125
$ws.Range("J135").Value = 1190.1
$ws.Range("L135").Value = 10710.9
$ws.Range("N135").Value = -15780.9
$ws.Range("H137").Value = 2042.375
$ws.Range("I137").Value = 768.2727
$ws.Range("J137").Value = 3120.4614
$ws.Range("K137").Value = 2304.8181
$ws.Range("L137").Value = 9361.3842
$ws.Range("M137").Value = 2795.1819
$ws.Range("N137").Value = -19561.3842

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 3391
$ws.Range("J9").Value = 5079
$ws.Range("L9").Value = 5079
$ws.Range("N9").Value = -5419
$ws.Range("H97").Value = 1311.7428
$ws.Range("I97").Value = 1250.2174
$ws.Range("K97").Value = 1250.2174
$ws.Range("M97").Value = -754.2174
$ws.Range("H102").Value = 962799.5600000001
$ws.Range("I102").Value = 1804199
$ws.Range("J102").Value = 6663.773
$ws.Range("K102").Value = 1804199
$ws.Range("L102").Value = 6663.773
$ws.Range("M102").Value = -1802577
$ws.Range("N102").Value = -9907.773000000001
$ws.Range("H112").Value = 74900
$ws.Range("J112").Value = 74900
$ws.Range("L112").Value = 74900
$ws.Range("N112").Value = -77116
$ws.Range("H113").Value = 4470.16
$ws.Range("I113").Value = 1460.7894
$ws.Range("J113").Value = 13999.833
$ws.Range("K113").Value = 1460.7894
$ws.Range("L113").Value = 13999.833
$ws.Range("M113").Value = 709.2106000000001
$ws.Range("N113").Value = -18339.833
$ws.Range("H126").Value = 11769804
$ws.Range("I126").Value = 14708550
$ws.Range("K126").Value = 44125650
$ws.Range("M126").Value = -44123180
$ws.Range("H127").Value = 15000
$ws.Range("J127").Value = 15000
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920
$ws.Range("H132").Value = 35718676
$ws.Range("I132").Value = 50004400
$ws.Range("K132").Value = 150013200
$ws.Range("M132").Value = -150010670

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 274.5
$ws.Range("I9").Value = 99
$ws.Range("K9").Value = 99
$ws.Range("M9").Value = 125
$ws.Range("H68").Value = 2629
$ws.Range("I68").Value = 2080.7693
$ws.Range("J68").Value = 3816.8333
$ws.Range("K68").Value = 2080.7693
$ws.Range("L68").Value = 3816.8333
$ws.Range("M68").Value = -1331.7693
$ws.Range("N68").Value = -5314.8333
$ws.Range("H71").Value = 2629
$ws.Range("I71").Value = 2080.7693
$ws.Range("J71").Value = 3816.8333
$ws.Range("K71").Value = 10403.8465
$ws.Range("L71").Value = 19084.1665
$ws.Range("M71").Value = -6659.8465
$ws.Range("N71").Value = -26572.1665
$ws.Range("H93").Value = 4523.364
$ws.Range("I93").Value = 4219.625
$ws.Range("J93").Value = 5333.3335
$ws.Range("K93").Value = 4219.625
$ws.Range("L93").Value = 5333.3335
$ws.Range("M93").Value = -2971.625
$ws.Range("N93").Value = -7829.3335
$ws.Range("H122").Value = 4593
$ws.Range("I122").Value = 4368.875
$ws.Range("J122").Value = 4849.143
$ws.Range("K122").Value = 13106.625
$ws.Range("L122").Value = 14547.429
$ws.Range("M122").Value = -10656.625
$ws.Range("N122").Value = -19447.429
$ws.Range("H132").Value = 4907.975
$ws.Range("I132").Value = 4917.5
$ws.Range("J132").Value = 4902.846
$ws.Range("K132").Value = 14752.5
$ws.Range("L132").Value = 14708.538
$ws.Range("M132").Value = -12222.5
$ws.Range("N132").Value = -19768.538
$ws.Range("H136").Value = 42864452
$ws.Range("I136").Value = 31257720
$ws.Range("K136").Value = 93773160
$ws.Range("M136").Value = -93770610

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3850596.8
$ws.Range("I23").Value = 5002380.5
$ws.Range("K23").Value = 5002380.5
$ws.Range("M23").Value = -5002151.5
$ws.Range("H107").Value = 6452434.5
$ws.Range("I107").Value = 10000573
$ws.Range("K107").Value = 30001719
$ws.Range("M107").Value = -29999799
$ws.Range("H122").Value = 12940.294
$ws.Range("I122").Value = 6049
$ws.Range("J122").Value = 22785
$ws.Range("K122").Value = 18147
$ws.Range("L122").Value = 68355
$ws.Range("M122").Value = -15697
$ws.Range("N122").Value = -73255
$ws.Range("H126").Value = 2912.4583
$ws.Range("I126").Value = 1743.2858
$ws.Range("J126").Value = 4549.3
$ws.Range("K126").Value = 5229.857400000001
$ws.Range("L126").Value = 13647.9
$ws.Range("M126").Value = -2759.857400000001
$ws.Range("N126").Value = -18587.9
$ws.Range("H132").Value = 5248.2
$ws.Range("I132").Value = 4727.0303
$ws.Range("K132").Value = 14181.0909
$ws.Range("M132").Value = -11651.0909
$ws.Range("H136").Value = 15165698
$ws.Range("I136").Value = 23821416
$ws.Range("J136").Value = 18191.416
$ws.Range("K136").Value = 71464248
$ws.Range("L136").Value = 54574.24800000001
$ws.Range("M136").Value = -71461698
$ws.Range("N136").Value = -59674.24800000001
